$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new column before D ("Hostname"), shifting the old
#        D..G (Zone, Gateway, OS, Extra) right into E..H. -----------------
$ws.Columns("D").Insert()

# --- 2. Header row -----------------------------------------------------
$ws.Cells.Item(1,4).Value = "Hostname"

# --- 3. Data rows 2-10 ---------------------------------------------------
# Columns: A=Hardware B=Port C=IPaddress D=Hostname E=Zone F=Gateway G=OS H=Extra
$rows = @(
  @{r=2;  c="10.1.205.11/16"; d=$null;        g="XenServer 7.3";  h="IMM"}
  @{r=3;  c="10.1.205.12/16"; d=$null;        g="XenServer 7.3";  h="XenServer Management"}
  @{r=4;  c="172.16.0.45/16"; d="ShareFile";  g="VIP NetScaler 12"; h="NetScaler DMZ"}
  @{r=5;  c="10.1.205.2/16";  d="netscaler";  g="VIP NetScaler 12"; h="Netscaler Management"}
  @{r=6;  c="10.1.205.4/16";  d="netscaler";  g="VIP NetScaler 12"; h="Netscaler Internal"}
  @{r=7;  c="10.1.205.1/16";  d="ldap";       g="Windows Server 2016 Standard"; h="Virutal LDAP server"}
  @{r=8;  c="10.1.205.3/16";  d="fileserver"; g="Windows Server 2016 Standard"; h="Virtual file server"}
  @{r=9;  c="10.1.205.10/16"; d=$null;        g="Citrix License Server Virtual Appliance v11.14.0.1_19800"; h="Virtual License Server"}
  @{r=10; c="10.1.205.5/16";  d=$null;        g=$null;            h="AAA Virtual Server"}
)

foreach ($row in $rows) {
  $r = $row.r
  $ws.Cells.Item($r,1).Value = "IBM x3650 M3 7945G2G"
  if ($r -eq 4) {
    $ws.Cells.Item($r,2).Value = 2
  } elseif ($r -eq 2) {
    $ws.Cells.Item($r,2).Value = 0
  } else {
    $ws.Cells.Item($r,2).Value = 1
  }
  $ws.Cells.Item($r,3).Value = $row.c
  if ($row.d -ne $null) {
    $ws.Cells.Item($r,4).Value = $row.d
  }
  $ws.Cells.Item($r,5).Value = "Intern"
  $ws.Cells.Item($r,6).Value = "10.1.5.32"
  if ($row.g -ne $null) {
    $ws.Cells.Item($r,7).Value = $row.g
  }
  $ws.Cells.Item($r,8).Value = $row.h
}

# Row 4 and 5 have Zone/Gateway overridden (DMZi zone) -------------------
$ws.Cells.Item(4,5).Value = "DMZi"
$ws.Cells.Item(4,6).Value = "172.16.0.35"

# --- 4. Column widths ------------------------------------------------
$ws.Columns("D").ColumnWidth = 15.140625
$ws.Columns("G").ColumnWidth = 54.7109375

# --- 5. New blank row 27 (matches rows 24-26 pattern) -------------------
$ws.Cells.Item(27,1).Value = "x"
$ws.Cells.Item(27,1).ClearContents()

# --- 6. Selection --------------------------------------------------------
$ws.Range("C10").Select()
